# Kruskal-Wallis, Dunn's test update
# Adds a new "Kruskal Wallis test" sheet (after STATS) containing a
# Kruskal-Wallis summary table and a Dunn's-test multiple-comparisons table.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet (STATS) ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Kruskal Wallis test"

# --- Pre-seed brand-new shared strings in the exact order they should
#     first appear so the shared-string table layout matches the target ---
$ws.Range("Z1").Value = "H-statistic"
$ws.Range("Z2").Value = "**"
$ws.Range("Z3").Value = "figure"
$ws.Range("Z4").Value = "Kruskal Wallis Test"
$ws.Range("Z5").Value = "Dunn's Test for Multiple Comparisons"
$ws.Range("Z6").Value = "Dunn's Test (No Correction)"

# ===================== Table 1: Kruskal Wallis Test =====================

$ws.Range("A1").Value = "Kruskal Wallis Test"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1:G1").Merge()

$ws.Range("B2").Value = "figure"
$ws.Range("C2").Value = "phase"
$ws.Range("D2").Value = "y-value"
$ws.Range("E2").Value = "H-statistic"
$ws.Range("F2").Value = "p-value"
$ws.Range("G2").Value = "significance"

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "SuppFig 2G"
$ws.Range("C3").Value = "interphase"
$ws.Range("D3").Value = "AC/Cyto RFP"
$ws.Range("E3").Value = 9.5737654948181401
$ws.Range("F3").Value = 0.0083384098404025303
$ws.Range("G3").Value = "**"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "SuppFig 2H"
$ws.Range("C4").Value = "interphase"
$ws.Range("D4").Value = "AC/BC RFP"
$ws.Range("E4").Value = 18.484149563096899
$ws.Range("F4").Value = 0.000096876384231799195
$ws.Range("F4").NumberFormat = "0.00E+00"
$ws.Range("G4").Value = "****"

# ============== Table 2: Dunn's Test for Multiple Comparisons ==============

$ws.Range("A6").Value = "Dunn's Test for Multiple Comparisons"
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6:R6").Merge()

$ws.Range("B7").Value = "figure"
$ws.Range("C7").Value = "control variable"
$ws.Range("D7").Value = "experimental variable"
$ws.Range("E7").Value = "phase"
$ws.Range("F7").Value = "test"
$ws.Range("G7").Value = "y-value"
$ws.Range("H7").Value = "control n"
$ws.Range("I7").Value = "experimental n"
$ws.Range("J7").Value = "control mean"
$ws.Range("K7").Value = "experimental mean"
$ws.Range("L7").Value = "control stdev"
$ws.Range("M7").Value = "experimental stdev"
$ws.Range("N7").Value = "test statistic"
$ws.Range("O7").Value = "degrees of freedom"
$ws.Range("P7").Value = "critical value"
$ws.Range("Q7").Value = "p-value"
$ws.Range("R7").Value = "significance"

$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "SuppFig 2G"
$ws.Range("C8").Value = "WT"
$ws.Range("D8").Value = "PP4r2 RNAi, 0X MTOC"
$ws.Range("E8").Value = "interphase"
$ws.Range("F8").Value = "Dunn's Test (No Correction)"
$ws.Range("G8").Value = "AC/Cyto RFP"
$ws.Range("H8").Value = 28
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 2.3574655660714199
$ws.Range("K8").Value = 1.2627672783333299
$ws.Range("L8").Value = 0.60595442502166297
$ws.Range("M8").Value = 0.65400877896671805
$ws.Range("N8").Value = "N/A"
$ws.Range("O8").Value = "N/A"
$ws.Range("P8").Value = "N/A"
$ws.Range("Q8").Value = 0.0021268200000000002
$ws.Range("Q8").NumberFormat = "0.00E+00"
$ws.Range("R8").Value = "**"

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "SuppFig 2G"
$ws.Range("C9").Value = "WT"
$ws.Range("D9").Value = "PP4r2 RNAi, 2X MTOC"
$ws.Range("E9").Value = "interphase"
$ws.Range("F9").Value = "Dunn's Test (No Correction)"
$ws.Range("G9").Value = "AC/Cyto RFP"
$ws.Range("H9").Value = 28
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 2.3574655660714199
$ws.Range("K9").Value = 2.3308286133333298
$ws.Range("L9").Value = 0.60595442502166297
$ws.Range("M9").Value = 0.19634361040205001
$ws.Range("N9").Value = "N/A"
$ws.Range("O9").Value = "N/A"
$ws.Range("P9").Value = "N/A"
$ws.Range("Q9").Value = 0.84113199999999999
$ws.Range("Q9").NumberFormat = "0.00E+00"
$ws.Range("R9").Value = "ns"

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "SuppFig 2H"
$ws.Range("C10").Value = "WT"
$ws.Range("D10").Value = "PP4r2 RNAi, 0X MTOC"
$ws.Range("E10").Value = "interphase"
$ws.Range("F10").Value = "Dunn's Test (No Correction)"
$ws.Range("G10").Value = "AC/BC RFP"
$ws.Range("H10").Value = 28
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 2.6261093049999999
$ws.Range("K10").Value = 1.1540538783333301
$ws.Range("L10").Value = 0.89850663775836404
$ws.Range("M10").Value = 0.24106337359836399
$ws.Range("N10").Value = "N/A"
$ws.Range("O10").Value = "N/A"
$ws.Range("P10").Value = "N/A"
$ws.Range("Q10").Value = 0.00035854300000000001
$ws.Range("Q10").NumberFormat = "0.00E+00"
$ws.Range("R10").Value = "***"

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "SuppFig 2H"
$ws.Range("C11").Value = "WT"
$ws.Range("D11").Value = "PP4r2 RNAi, 2X MTOC"
$ws.Range("E11").Value = "interphase"
$ws.Range("F11").Value = "Dunn's Test (No Correction)"
$ws.Range("G11").Value = "AC/BC RFP"
$ws.Range("H11").Value = 28
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 2.6261093049999999
$ws.Range("K11").Value = 0.95817877333333301
$ws.Range("L11").Value = 0.89850663775836404
$ws.Range("M11").Value = 0.155808714255736
$ws.Range("N11").Value = "N/A"
$ws.Range("O11").Value = "N/A"
$ws.Range("P11").Value = "N/A"
$ws.Range("Q11").Value = 0.00501229
$ws.Range("Q11").NumberFormat = "0.00E+00"
$ws.Range("R11").Value = "**"

# --- Clear the scratch cells used to fix shared-string insertion order ---
$ws.Range("Z1:Z6").Clear()

# --- Selection matches the freshly-added sheet's final cursor position ---
$ws.Range("B13").Select()

Write-Output "done"
